$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update Owner column (E) for specific rows
$ws.Range("E8").Value  = "Chief Technology Officer"
$ws.Range("E9").Value  = "IT Managers"
$ws.Range("E10").Value = "DevOps Engineers"
$ws.Range("E11").Value = "System Administrators"
$ws.Range("E14").Value = "Chief Technology Officer"
$ws.Range("E15").Value = "IT Managers"
$ws.Range("E16").Value = "DevOps Engineers"
$ws.Range("E17").Value = "System Administrators"

# Update Dependencies (I) and Notes (J) columns for rows 8 through 17
for ($row = 8; $row -le 17; $row++) {
    $ws.Range("I$row").Value = "Dependent on Cloud Infrastructure Migration milestone completion"
    $ws.Range("J$row").Value = "Critical action for Information Technology success"
}
